# Updates Price (D) and Volume(1h) (E) columns with refreshed crypto figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.717.53"
$ws.Range("E2").Value = "  -1.67%  "

$ws.Range("D3").Value = "3.638.10"
$ws.Range("E3").Value = "  -2.01%  "

$ws.Range("E4").Value = "  +18.93%  "

$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "225.19"
$ws.Range("E6").Value = "  -5.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "646.45"
$ws.Range("E7").Value = "  -1.49%  "

$ws.Range("E8").Value = "  -3.95%  "

$ws.Range("E9").Value = "  +4.38%  "

$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("D11").Value = "3.635.41"
$ws.Range("E11").Value = "  -2.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.00"
$ws.Range("E12").Value = "  +14.82%  "

$ws.Range("E13").Value = "  +5.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000298"
$ws.Range("E14").Value = "  -3.70%  "

$ws.Range("E15").Value = "  -4.01%  "

$ws.Range("D16").Value = "4.317.71"
$ws.Range("E16").Value = "  -2.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.89"
$ws.Range("E17").Value = "  +31.28%  "

$ws.Range("D18").Value = "95.510.90"
$ws.Range("E18").Value = "  -1.70%  "

$ws.Range("E19").Value = "  -5.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.92"
$ws.Range("E20").Value = "  +6.42%  "

$ws.Range("D21").Value = "3.626.16"
$ws.Range("E21").Value = "  -2.47%  "

$ws.Range("E22").Value = "  +49.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.540"
$ws.Range("E23").Value = "  -2.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "536.78"
$ws.Range("E24").Value = "  +1.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.01"
$ws.Range("E25").Value = "  +12.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.38"
$ws.Range("E26").Value = "  -2.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.22"
$ws.Range("E27").Value = "  +4.40%  "

$ws.Range("E28").Value = "  -8.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.73"
$ws.Range("E29").Value = "  +1.87%  "

$ws.Range("D30").Value = "3.807.08"
$ws.Range("E30").Value = "  -2.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.73"
$ws.Range("E31").Value = "  +6.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  +5.75%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.93"
$ws.Range("E34").Value = "  +5.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.645"
$ws.Range("E35").Value = "  +6.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "33.82"
$ws.Range("E36").Value = "  +1.98%  "

$ws.Range("E37").Value = "  -4.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.22%  "

$ws.Range("E39").Value = "  +22.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.65"
$ws.Range("E40").Value = "  -1.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "605.65"
$ws.Range("E41").Value = "  -5.37%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.26"
$ws.Range("E43").Value = "  +5.20%  "

$ws.Range("E44").Value = "  +3.65%  "

$ws.Range("E45").Value = "  +6.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.55"
$ws.Range("E46").Value = "  +1.95%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.05"
$ws.Range("E47").Value = "  +1.19%  "

$ws.Range("E48").Value = "  -6.27%  "

$ws.Range("E49").Value = "  +4.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "237.58"
$ws.Range("E50").Value = "  +13.37%  "

$ws.Range("E51").Value = "  -0.61%  "
